$wb = $excel.ActiveWorkbook

# Update the handoff/handback timestamps for the 94c9d3cb-... localization
# entry (row 3) on each language sheet, reflecting a newly generated
# handback report.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-10 00:54:13"
$wsZhCn.Range("G3").Value = "2016-03-10 00:55:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-10 00:54:19"
$wsDeDe.Range("G3").Value = "2016-03-10 00:55:33"
